$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table of lines/contingencies grows from 14 rows (line1-6, extr1-8) to 16 rows
# (line1-8, extr1-8): two new "line7"/"line8" records are inserted logically before
# extr1, pushing the extr* data down by two rows, and a couple of values are refined
# ("rene fine") along the way.

# First, extend new rows 16 and 17 with the same formatting (bold border style) used
# by the rest of column A, so no new/duplicate style gets introduced.
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(17, 1))

$data = @(
    @(0, "line1", 7, 9, $true),
    @(1, "line2", 9, 8, $true),
    @(2, "line3", 8, 10, $true),
    @(3, "line4", 8, 11, $true),
    @(4, "line5", 10, 5, $true),
    @(5, "line6", 12, 8, $true),
    @(6, "line7", 14, 11, $true),
    @(7, "line8", 16, 9, $false),
    @(8, "extr1", 5, 12, $true),
    @(9, "extr2", 5, 9, $true),
    @(10, "extr3", 10, 11, $false),
    @(11, "extr4", 7, 8, $true),
    @(12, "extr5", 9, 11, $false),
    @(13, "extr6", 7, 11, $false),
    @(14, "extr7", 5, 7, $true),
    @(15, "extr8", 8, 5, $false)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
